$wb = $excel.ActiveWorkbook

# --- Rename Sheet3 -> PIToArduino -------------------------------------
$wsPI = $wb.Worksheets.Item("Sheet3")
$wsPI.Name = "PIToArduino"

$wsNet = $wb.Worksheets.Item("NetworkTables")
$wsArd = $wb.Worksheets.Item("FromArduinoToPi")

# --- NetworkTables: append the new "toPi" message block (rows 62-68) --
$wsNet.Range("B62").Value = "toPi"
$wsNet.Range("C62").Value = "LED"

$wsNet.Range("B63").Value = "toPi"
$wsNet.Range("D63").Value = "sensorLevel"
$wsNet.Range("C63").Value = "prox01"

$wsNet.Range("B64").Value = "toPi"
$wsNet.Range("D64").Value = "sensorLevel"
$wsNet.Range("C64").Value = "prox02"

$wsNet.Range("B65").Value = "toPi"
$wsNet.Range("D65").Value = "sensorLevel"
$wsNet.Range("C65").Value = "prox03"

$wsNet.Range("B66").Value = "toPi"
$wsNet.Range("D66").Value = "sensorLevel"
$wsNet.Range("C66").Value = "prox04"

$wsNet.Range("B67").Value = "toPi"
$wsNet.Range("D67").Value = "sensorLevel"
$wsNet.Range("C67").Value = "prox05"

$wsNet.Range("B68").Value = "toPi"
$wsNet.Range("D68").Value = "sensorLevel"
$wsNet.Range("C68").Value = "prox06"

# --- Update each sheet's view/selection state --------------------------
$wsArd.Select()
$wsArd.Range("C4:D16").Select()

$wsPI.Select()
$wsPI.Range("C3:D9").Select()

# Leave NetworkTables as the active/selected sheet, positioned on the
# newly-added block, matching the authored selection.
$wsNet.Select()
$wsNet.Range("C69").Select()
